$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "index" cell (A2) so the new rows'
# A-cells pick up the same bold/centered/bordered style (s="1") without
# minting new style records.
$ws.Range("A2").Copy()

# Row 3 - "linear"
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "linear"
$ws.Range("C3").Value = 260.57
$ws.Range("D3").Value = 121449.22
$ws.Range("E3").Value = 0.64
$ws.Range("F3").Value = 202.88
$ws.Range("G3").Value = 79043.21000000001
$ws.Range("H3").Value = 0.73
$ws.Range("I3").Value = 136.06
$ws.Range("J3").Value = 29849.67
$ws.Range("K3").Value = 0.89

# Row 4 - "baseline-rent"
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "baseline-rent"
$ws.Range("C4").Value = 285.56
$ws.Range("D4").Value = 140135.42
$ws.Range("E4").Value = 0.58
$ws.Range("F4").Value = 237.83
$ws.Range("G4").Value = 92491.55
$ws.Range("H4").Value = 0.6899999999999999

$excel.CutCopyMode = $false
